$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Host / Username / Password columns (D, F, G) for data rows 2-7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 4).Value = "ipfdb.database.windows.net"
    $ws.Cells.Item($r, 6).Value = "ipfadmin"
    $ws.Cells.Item($r, 7).Value = "ifp@admin1"
}

# Add hyperlinks to the Password column (G2:G7), matching the mailto-style
# convention already used for BranchMail/HelpDeskMail columns, then restore
# the Hyperlink cell style so it matches the existing linked cells.
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $ws.Hyperlinks.Add($cell, "mailto:ifp@admin1")
    $cell.Style = "Hyperlink"
}

# Update the current selection
$ws.Range("F7:G7").Select()
